{"js": "// Replace each old three-digit-division expression with its new value.\nconst pairs = [\n  [\"499\u00f76=83, 1\", \"515\u00f79=57, 2\"],\n  [\"392\u00f79=43, 5\", \"910\u00f74=227, 2\"],\n  [\"625\u00f78=78, 1\", \"558\u00f73=186, 0\"],\n  [\"367\u00f74=91, 3\", \"634\u00f72=317, 0\"],\n  [\"765\u00f72=382, 1\", \"973\u00f73=324, 1\"],\n  [\"819\u00f77=117, 0\", \"697\u00f73=232, 1\"],\n  [\"138\u00f72=69, 0\", \"918\u00f75=183, 3\"],\n  [\"585\u00f76=97, 3\", \"312\u00f74=78, 0\"],\n  [\"875\u00f79=97, 2\", \"875\u00f73=291, 2\"],\n  [\"634\u00f75=126, 4\", \"174\u00f72=87, 0\"],\n  [\"374\u00f78=46, 6\", \"463\u00f76=77, 1\"],\n  [\"289\u00f73=96, 1\", \"350\u00f73=116, 2\"],\n  [\"259\u00f78=32, 3\", \"578\u00f74=144, 2\"],\n  [\"721\u00f77=103, 0\", \"421\u00f72=210, 1\"],\n  [\"161\u00f76=26, 5\", \"732\u00f75=146, 2\"],\n  [\"949\u00f75=189, 4\", \"678\u00f76=113, 0\"],\n  [\"822\u00f76=137, 0\", \"218\u00f72=109, 0\"],\n  [\"737\u00f79=81, 8\", \"639\u00f72=319, 1\"],\n  [\"960\u00f76=160, 0\", \"306\u00f76=51, 0\"],\n  [\"409\u00f79=45, 4\", \"827\u00f76=137, 5\"],\n  [\"996\u00f79=110, 6\", \"214\u00f79=23, 7\"],\n  [\"895\u00f77=127, 6\", \"878\u00f72=439, 0\"],\n  [\"286\u00f77=40, 6\", \"401\u00f73=133, 2\"],\n  [\"342\u00f79=38, 0\", \"823\u00f76=137, 1\"],\n  [\"908\u00f73=302, 2\", \"430\u00f76=71, 4\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('499\u00f76=83, 1', '515\u00f79=57, 2'),\n    @('392\u00f79=43, 5', '910\u00f74=227, 2'),\n    @('625\u00f78=78, 1', '558\u00f73=186, 0'),\n    @('367\u00f74=91, 3', '634\u00f72=317, 0'),\n    @('765\u00f72=382, 1', '973\u00f73=324, 1'),\n    @('819\u00f77=117, 0', '697\u00f73=232, 1'),\n    @('138\u00f72=69, 0', '918\u00f75=183, 3'),\n    @('585\u00f76=97, 3', '312\u00f74=78, 0'),\n    @('875\u00f79=97, 2', '875\u00f73=291, 2'),\n    @('634\u00f75=126, 4', '174\u00f72=87, 0'),\n    @('374\u00f78=46, 6', '463\u00f76=77, 1'),\n    @('289\u00f73=96, 1', '350\u00f73=116, 2'),\n    @('259\u00f78=32, 3', '578\u00f74=144, 2'),\n    @('721\u00f77=103, 0', '421\u00f72=210, 1'),\n    @('161\u00f76=26, 5', '732\u00f75=146, 2'),\n    @('949\u00f75=189, 4', '678\u00f76=113, 0'),\n    @('822\u00f76=137, 0', '218\u00f72=109, 0'),\n    @('737\u00f79=81, 8', '639\u00f72=319, 1'),\n    @('960\u00f76=160, 0', '306\u00f76=51, 0'),\n    @('409\u00f79=45, 4', '827\u00f76=137, 5'),\n    @('996\u00f79=110, 6', '214\u00f79=23, 7'),\n    @('895\u00f77=127, 6', '878\u00f72=439, 0'),\n    @('286\u00f77=40, 6', '401\u00f73=133, 2'),\n    @('342\u00f79=38, 0', '823\u00f76=137, 1'),\n    @('908\u00f73=302, 2', '430\u00f76=71, 4'),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
